$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.740322580645161
$ws.Range("C2").Value = 0.718626155878468
$ws.Range("D2").Value = 0.753222836095764
$ws.Range("E2").Value = 0.708074534161491
$ws.Range("F2").Value = 0.695043511161559

$ws.Range("B3").Value = 0.614516129032258
$ws.Range("C3").Value = 0.730515191545575
$ws.Range("D3").Value = 0.692449355432781
$ws.Range("E3").Value = 0.68167701863354
$ws.Range("F3").Value = 0.617858494135452

$ws.Range("B4").Value = 0.620967741935484
$ws.Range("C4").Value = 0.624834874504623
$ws.Range("D4").Value = 0.591160220994475
$ws.Range("E4").Value = 0.596273291925466
$ws.Range("F4").Value = 0.540295119182747

$ws.Range("B5").Value = 0.885483870967742
$ws.Range("C5").Value = 0.919418758256275
$ws.Range("D5").Value = 0.896869244935543
$ws.Range("E5").Value = 0.875776397515528
$ws.Range("F5").Value = 0.865304578130912
